$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure updated cells keep their original "Text" storage type so values
# such as "278.99" or "1.08%" are not reinterpreted as numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.08%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.14%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.845"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.28%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06384"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.50%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.033"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.301"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.97%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8931"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.88%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.52%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05785"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.38%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07507"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.15%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02918"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.60%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08998"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.26%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001588"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.53%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006381"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.42%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006047"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.73%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.474"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.54%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.305"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.08%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.303"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.35%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1348"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.43%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.899"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.06%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1506"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "9.10%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04385"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.48%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001175"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.30%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "10.76%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001179"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.67%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-14.59%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04046"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.99%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1414"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "19.97%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006650"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.37%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002058"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.98%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01115"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.22%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005560"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.56%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.628"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.53%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-19.57%"
